$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Models")

# Reshape the Models summary table: columns B..F previously held
# [Command, r2, ar2, f-complexity, q-complexity]; the sheet now puts the
# numeric stats first (r2, ar2, f-complexity, q-complexity) and moves the
# Command string out to column F.

$ws.Range("B1").Value = "r2"
$ws.Range("C1").Value = "ar2"
$ws.Range("D1").Value = "f-complexity"
$ws.Range("E1").Value = "q-complexity"
$ws.Range("F1").Value = "Command"
$ws.Range("B2").ClearContents()
$ws.Range("D2").Value = 59
$ws.Range("E2").Value = 10
$ws.Range("F2").Value = "reg q2 eq* boughtSample employer male unemployed _region* _income* _stem* _industry* _age* cage* cincome* cprovider*"
$ws.Range("B6").Value = 0.42
$ws.Range("C6").Value = -0.09
$ws.Range("D6").Value = 58
$ws.Range("E6").Value = 10
$ws.Range("F6").Value = "reg q2 eq* boughtSample employer male unemployed _region* _income* _stem* _industry* _age* cage* cincome* cprovider*"
$ws.Range("B7").Value = 0.39
$ws.Range("C7").Value = 0.2
$ws.Range("D7").Value = 23
$ws.Range("E7").Value = 10
$ws.Range("F7").Value = "reg q2 eq4 eq6squared eq6cubed male unemployed _region3 _region5 _region6 _income2 _income3 _income6 _stem1 _stem2 _industry2 _industry6 _industry7 _industry9 _industry12 _age2 cage2 cincome3 cprovider1"
$ws.Range("B8").Value = 0.37
$ws.Range("C8").Value = 0.21
$ws.Range("D8").Value = 20
$ws.Range("E8").Value = 9
$ws.Range("F8").Value = "reg q2 eq4 eq6squared eq6cubed unemployed _region3 _region6 _income2 _income3 _income6 _stem1 _stem2 _industry2 _industry6 _industry7 _industry9 _industry12 _age2 cincome3 cprovider1"
$ws.Range("B9").Value = 0.18
$ws.Range("C9").Value = 0.13
$ws.Range("D9").Value = 6
$ws.Range("E9").Value = 3
$ws.Range("F9").Value = "reg q2 eq6squared eq6cubed _income2 _industry6 _industry9"
$ws.Range("B10").Value = 0.49
$ws.Range("C10").Value = 0.02
$ws.Range("D10").Value = 59
$ws.Range("E10").Value = 10
$ws.Range("F10").Value = "reg index eq* boughtSample employer male unemployed _region* _income* _stem* _industry* _age* cage* cincome* cprovider*"
$ws.Range("B11").Value = 0.48
$ws.Range("C11").Value = 0.25
$ws.Range("D11").Value = 30
$ws.Range("E11").Value = 10
$ws.Range("F11").Value = "reg index eq4 eq6 eq6cubed employer male unemployed _region3 _region5 _region6 _region8 _income2 _income6 _income10 _stem1 _stem2 _industry2 _industry4-_industry6 _industry9-_industry12 cage2 cage3 cincome2 cincome3 cprovider1 cprovider2"
$ws.Range("B12").Value = 0.44
$ws.Range("C12").Value = 0.3
$ws.Range("D12").Value = 20
$ws.Range("E12").Value = 10
$ws.Range("F12").Value = "reg index eq4 eq6 eq6cubed male unemployed _region3 _region5 _region6 _region8 _income6 _stem2 _industry6 _industry9 _industry11 cage2 cage3 cincome2 cincome3 cprovider1"
$ws.Range("B13").Value = 0.27
$ws.Range("C13").Value = 0.23
$ws.Range("D13").Value = 7
$ws.Range("E13").Value = 4
$ws.Range("F13").Value = "reg index eq4 eq6 eq6cubed _income6 _industry6 _industry9"
$ws.Range("B14").Value = 0.44
$ws.Range("C14").Value = -0.05
$ws.Range("D14").Value = 58
$ws.Range("E14").Value = 10
$ws.Range("F14").Value = "reg index eq* boughtSample employer male unemployed _region* _income* _stem* _industry* _age* cage* cincome* cprovider*"
$ws.Range("B15").Value = 0.39
$ws.Range("C15").Value = 0.21
$ws.Range("D15").Value = 23
$ws.Range("E15").Value = 10
$ws.Range("F15").Value = "reg index eq4squared eq4cubed eq6 eq6cubed employer male unemployed _region3 _region5 _income9 _stem2 _industry2 _industry4 _industry5 _industry6  _industry9- _industry12 _age2 cprovider1 cprovider2"
$ws.Range("B16").Value = 0.35
$ws.Range("C16").Value = 0.25
$ws.Range("D16").Value = 14
$ws.Range("E16").Value = 7
$ws.Range("F16").Value = "reg index eq4squared eq4cubed eq6 eq6cubed male unemployed _region3 _region5 _industry4 _industry6  _industry9 _industry12 cprovider1"
$ws.Range("B17").Value = 0.22
$ws.Range("C17").Value = 0.18
$ws.Range("D17").Value = 7
$ws.Range("E17").Value = 3
$ws.Range("F17").Value = "reg index eq4squared eq4cubed eq6 eq6cubed _industry6  _industry9"
$ws.Range("B18").Value = 0.34
$ws.Range("C18").Value = -0.12
$ws.Range("D18").Value = 53
$ws.Range("E18").Value = 8
$ws.Range("F18").Value = "reg index boughtSample employer male unemployed _region* _income* _stem* _industry* _age* cage* cincome* cprovider*"
$ws.Range("B19").Value = 0.32
$ws.Range("C19").Value = 0.16
$ws.Range("D19").Value = 19
$ws.Range("E19").Value = 8
$ws.Range("F19").Value = "reg index male unemployed _region8 _region9 _income6 _income8 _stem1 _stem2 _industry2 _industry4-_industry6 _industry9 _industry11 cage2 cage3 cprovider1 cprovider2"
$ws.Range("B20").Value = 0.28
$ws.Range("C20").Value = 0.18
$ws.Range("D20").Value = 13
$ws.Range("E20").Value = 8
$ws.Range("F20").Value = "reg index male unemployed _region8 _income6 _income8 _stem1 _stem2 _industry6 _industry9 cage2 cage3 cprovider1"
$ws.Range("B21").Value = 0.17
$ws.Range("C21").Value = 0.13
$ws.Range("D21").Value = 4
$ws.Range("E21").Value = 4
$ws.Range("F21").Value = "reg index male _income6 _industry9 cprovider1"

# Column B no longer carries the long Command text, so its explicit
# custom width is no longer meaningful there; column F (which now holds
# Command) gets the width that used to live on column B.
$ws.Columns.Item(6).ColumnWidth = $ws.Columns.Item(2).ColumnWidth

# Move the cursor/selection to D3, matching the author's last-touched cell.
$ws.Range("D3").Select()
